$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Matrix - result")
$ws2.Select()
$r = $excel.ActiveWindow.ScrollColumn
$c = $excel.ActiveWindow.ScrollRow
